$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply the autofilter so it spans the (pre-existing) data range A1:B54 ---
# (done before adding the new row below it, otherwise the filter range would
#  auto-expand to include the newly added adjacent row)
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:B54").AutoFilter()

# --- Add the new "Electric Elemental" species row (A55/B55) ---
$ws.Range("A55").Value2 = "Electric Elemental"
$ws.Range("B55").Value2 = "Within the \imp{Elemental Planes} there can be found a single, enormous mountain, surrounded at all times by a roiling, black cloud filled with crackling energy: \key{Thundertop}. Lightning and thunder are everpresent in this hostile environment, and every surface is highly charged with static electricity – the foolish explorer who sets foot on the mountain of thunder without some rubber-soled boots is liable to have a {\it very} bad time. `nWithin the crackling chaos and the booming crashes of this formiddable environment, reside a number of creatures who have learned to harness, channel and consume electrical energy, using it for their own end - \key{Electric Elementals}. "

# Row height to fit the new description text
$ws.Rows.Item(55).RowHeight = 69

# --- Swap the two _xlnm._FilterDatabase defined names so the hidden one
#     (used internally by the filter) stores the full data range, and the
#     visible one reverts to just the header row ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        if ($n.Visible -eq $false) {
            $n.RefersTo = "=Sheet1!`$A`$1:`$B`$54"
        } else {
            $n.RefersTo = "=Sheet1!`$A`$1:`$B`$1"
        }
    }
}

# --- Move the selection/scroll position to the newly added row ---
[void]$ws.Range("B55").Select()
$excel.ActiveWindow.ScrollRow = 46
